$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.717.89"
$ws.Range("E2").Value = "  -1.51%  "
$ws.Range("D3").Value = "1.805.73"
$ws.Range("E3").Value = "  -1.03%  "
$ws.Range("E4").Value = "  +0.42%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "232.50"
$ws.Range("E5").Value = "  -0.98%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5930"
$ws.Range("E6").Value = "  -1.30%  "
$ws.Range("E7").Value = "  +0.43%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2776"
$ws.Range("E8").Value = "  -0.54%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06824"
$ws.Range("E9").Value = "  -3.36%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "23.35"
$ws.Range("E10").Value = "  -0.60%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07516"
$ws.Range("E11").Value = "  -1.15%  "
$ws.Range("D12").Value = "1.807.96"
$ws.Range("E12").Value = "  -1.03%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.758"
$ws.Range("E13").Value = "  -0.65%  "
$ws.Range("E14").Value = "  -0.97%  "
$ws.Range("D15").Value = "2.050.58"
$ws.Range("E15").Value = "  -1.07%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.000009228"
$ws.Range("E16").Value = "  -6.87%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "75.68"
$ws.Range("E17").Value = "  -3.97%  "
$ws.Range("D18").Value = "28.677.03"
$ws.Range("E18").Value = "  -1.66%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.476"
$ws.Range("E19").Value = "  -6.41%  "
$ws.Range("E20").Value = "  +0.37%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "211.13"
$ws.Range("E21").Value = "  -6.63%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "11.51"
$ws.Range("E22").Value = "  -1.71%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.828"
$ws.Range("E23").Value = "  -2.33%  "
$ws.Range("E24").Value = "  +0.42%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "154.00"
$ws.Range("E25").Value = "  -0.58%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.877"
$ws.Range("E26").Value = "  -1.72%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1267"
$ws.Range("E27").Value = "  -2.48%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "16.43"
$ws.Range("E28").Value = "  -0.72%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.428"
$ws.Range("E29").Value = "  -4.12%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.06163"
$ws.Range("E30").Value = "  -1.42%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.429"
$ws.Range("E31").Value = "  -1.42%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.789"
$ws.Range("E32").Value = "  -1.04%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.745"
$ws.Range("E33").Value = "  -1.31%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.734"
$ws.Range("E34").Value = "  -0.10%  "
$ws.Range("E35").Value = "  -5.29%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6422"
$ws.Range("E36").Value = "  +0.45%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.500"
$ws.Range("E37").Value = "  -1.23%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.717"
$ws.Range("E38").Value = "  -0.19%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.567"
$ws.Range("E39").Value = "  +1.06%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01696"
$ws.Range("E40").Value = "  -2.05%  "
$ws.Range("D41").Value = "1.145.56"
$ws.Range("E41").Value = "  -5.79%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8830"
$ws.Range("E42").Value = "  -2.23%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.007"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "99.92"
$ws.Range("E44").Value = "  -0.37%  "
$ws.Range("D45").Value = "1.957.10"
$ws.Range("E45").Value = "  -1.35%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "60.60"
$ws.Range("E46").Value = "  -3.35%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000112"
$ws.Range("E47").Value = "  -2.70%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.596"
$ws.Range("E48").Value = "  +0.04%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.370"
$ws.Range("E49").Value = "  -1.77%  "
$ws.Range("E50").Value = "  -0.44%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4478"
$ws.Range("E51").Value = "  -1.59%  "
